# "add data 23 Marzo" -- append the 23-March-2020 rows (one per
# Andalusian province) to the "datos" sheet and introduce a new
# "Curados" (recovered) column.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New column header (H1) -- this mints the new shared string "Curados".
$ws.Range("H1").Value = "Curados"

# New daily rows for 2020-03-23 (Excel serial date 43913), one row per
# province, in the same order used throughout the sheet.
$newRows = @(
    @{ Row = 90; Prov = "Almería"; C = 74;  D = 23;  E = 49;  F = 2;  G = 3;  H = 2 },
    @{ Row = 91; Prov = "Cádiz";   C = 178; D = 50;  E = 125; F = 3;  G = 6;  H = 4 },
    @{ Row = 92; Prov = "Córdoba"; C = 191; D = 50;  E = 137; F = 4;  G = 7;  H = 0 },
    @{ Row = 93; Prov = "Granada"; C = 374; D = 169; E = 188; F = 17; G = 18; H = 0 },
    @{ Row = 94; Prov = "Huelva";  C = 58;  D = 29;  E = 28;  F = 1;  G = 1;  H = 0 },
    @{ Row = 95; Prov = "Jaén";    C = 215; D = 79;  E = 131; F = 5;  G = 12; H = 5 },
    @{ Row = 96; Prov = "Málaga";  C = 520; D = 192; E = 307; F = 21; G = 44; H = 34 },
    @{ Row = 97; Prov = "Sevilla"; C = 351; D = 152; E = 194; F = 5;  G = 13; H = 5 }
)

foreach ($r in $newRows) {
    $i = $r.Row
    $ws.Range("A$i").Value = 43913
    $ws.Range("B$i").Value = $r.Prov
    $ws.Range("C$i").Value = $r.C
    $ws.Range("D$i").Value = $r.D
    $ws.Range("E$i").Value = $r.E
    $ws.Range("F$i").Value = $r.F
    $ws.Range("G$i").Value = $r.G
    $ws.Range("H$i").Value = $r.H
}

# Restore the view's active selection/scroll position to where the author
# left off editing.
[void]$ws.Range("I91").Select()
$excel.ActiveWindow.ScrollRow = 58
$excel.ActiveWindow.ScrollColumn = 1
